# Update "想去人数" (want-to-go count, column F) values on the "展览"
# sheet and the mirrored "全部类型" sheet, per the regenerated site data
# (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$exhibitionSheet = $wb.Worksheets.Item("展览")
$exhibitionSheet.Cells.Item(2, 6).Value  = 202
$exhibitionSheet.Cells.Item(4, 6).Value  = 5259
$exhibitionSheet.Cells.Item(13, 6).Value = 4202
$exhibitionSheet.Cells.Item(18, 6).Value = 3287
$exhibitionSheet.Cells.Item(20, 6).Value = 1073
$exhibitionSheet.Cells.Item(33, 6).Value = 21

$allTypesSheet = $wb.Worksheets.Item("全部类型")
$allTypesSheet.Cells.Item(2, 6).Value  = 202
$allTypesSheet.Cells.Item(5, 6).Value  = 5259
$allTypesSheet.Cells.Item(14, 6).Value = 4202
$allTypesSheet.Cells.Item(19, 6).Value = 3287
$allTypesSheet.Cells.Item(21, 6).Value = 1073
$allTypesSheet.Cells.Item(34, 6).Value = 21
